$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new values to the second (bottom) results table
$ws.Range("E9").Value = 0.26640588370690599
$ws.Range("F9").Value = 0.33193000603776401
$ws.Range("D10").Value = 0.18620018243242101
$ws.Range("E11").Value = 0.96447150925630099

# Update the active cell selection to match the author's final cursor position
$ws.Range("E11").Select()
